$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  2 = 0
  3 = 1
  4 = 2
  5 = 2
  6 = 0
  7 = 1
  8 = 2
  9 = 0
  10 = 5
  11 = 1
  12 = 1
  13 = 0
  14 = 1
  15 = 0
  16 = 2
  17 = 1
  18 = 0
  19 = 1
  20 = 0
  21 = 0
  23 = 1
  24 = 1
  25 = 0
  26 = 0
  27 = 2
  28 = 1
  29 = 1
  30 = 0
  31 = 2
  32 = 3
  33 = 0
  34 = 1
  35 = 1
  36 = 1
  37 = 1
  38 = 2
  39 = 0
  41 = 0
  42 = 1
  43 = 1
  44 = 2
  45 = 2
  46 = 0
  47 = 0
  48 = 1
  49 = 0
  50 = 1
  51 = 0
  52 = 1
  53 = 0
  54 = 0
  55 = 1
  56 = 2
  57 = 0
  58 = 2
  59 = 0
  60 = 1
  61 = 2
  62 = 3
  63 = 3
  64 = 1
  65 = 0
  66 = 0
  67 = 1
  68 = 2
  69 = 1
  70 = 2
  71 = 4
  72 = 0
  75 = 1
  76 = 2
  77 = 2
  78 = 2
  79 = 1
}

foreach ($row in $values.Keys) {
  $ws.Cells.Item($row, 7).Value = $values[$row]
}
